$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 261, pushing existing rows 261-333 down to 262-334
$ws.Rows("261:261").Insert()

# Populate the newly inserted row 261 with the new weekly record.
# Static/unchanged columns mirror the surrounding rows of this sub-dataset.
$ws.Range("A261").Value = 10
$ws.Range("B261").Value = "Vega Modelo de Temuco"
$ws.Range("C261").Value = "La Araucanía"
$ws.Range("D261").Value = 44722
$ws.Range("E261").Value = 9
$ws.Range("F261").Value = 100112009
$ws.Range("G261").Value = "Acelga"
$ws.Range("H261").Value = "Sin especificar"
$ws.Range("I261").Value = "Primera"
$ws.Range("J261").Value = 30
$ws.Range("K261").Value = 10000
$ws.Range("L261").Value = 10000
$ws.Range("M261").Value = 10000
$ws.Range("N261").Value = "`$/docena de atados (12 kilos)"
$ws.Range("O261").Value = "Provincia de Cautín"
$ws.Range("P261").Value = 833
$ws.Range("Q261").Value = 12
$ws.Range("R261").Value = "Hortaliza"

# Make sure the D261 cell keeps the date-style formatting (s="2") used by the
# rest of the Fecha column.
$ws.Range("D261").NumberFormat = $ws.Range("D262").NumberFormat
